$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.988.26'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '2.937.98'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '374.79'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').Value = '101.86'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('D7').Value = '0.535'
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('D10').Value = '36.49'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Value = '0.0837'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '3.396.39'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '17.92'
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '7.34'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('D16').Value = '2.939.82'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '0.975'
$ws.Range('E17').Value = '  +2.69%  '
$ws.Range('D18').Value = '50.907.43'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('E19').Value = '  -5.16%  '
$ws.Range('D20').Value = '7.16'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '264.18'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('D26').Value = '8.20'
$ws.Range('E26').Value = '  +10.98%  '
$ws.Range('D27').Value = '7.76'
$ws.Range('E27').Value = '  +9.08%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').Value = '25.61'
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('D32').Value = '9.84'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').Value = '50.77'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = '33.55'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0449'
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('E36').Value = '  -3.02%  '
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').Value = '16.36'
$ws.Range('E41').Value = '  -3.59%  '
$ws.Range('D42').Value = '1.79'
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').Value = '120.25'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').Value = '0.289'
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('D45').Value = '20.90'
$ws.Range('E45').Value = '  -4.52%  '
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '3.25'
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '2.31'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('D49').Value = '1.972.22'
$ws.Range('E49').Value = '  -2.81%  '
$ws.Range('D50').Value = '0.0342'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E51').Value = '  -0.87%  '
